$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Coin / Link / Volume(1h) cells (plain text, safe to assign directly) ---
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").Value = "  +3.61%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("E12").Value = "  +6.97%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").Value = "  +4.92%  "
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("E18").Value = "  +4.56%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("E23").Value = "  +35.24%  "
$ws.Range("E24").Value = "  +7.84%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("E27").Value = "  +12.52%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("E37").Value = "  +6.94%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("E46").Value = "  +7.01%  "
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("E48").Value = "  +7.69%  "
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("E50").Value = "  +6.84%  "
$ws.Range("E51").Value = "  +3.74%  "

# --- Update Price cells (column D) ---
# Several of the new prices are plain numeric-looking strings (e.g. "1.001").
# Assigning such a string straight to .Value lets Excel auto-convert it into a
# floating point number, which would lose the original text formatting used by
# this sheet (prices are stored as literal text, not numbers). To preserve the
# text representation exactly, write each value as a `="..."` text formula and
# then use Copy + PasteSpecial (values only) to collapse it back down to a plain
# literal string cell - this avoids Excel's automatic text->number coercion and
# does not require touching any cell/number formats.
$ws.Range("D2").Formula = '="30.733.58"'
$ws.Range("D3").Formula = '="1.894.18"'
$ws.Range("D4").Formula = '="1.001"'
$ws.Range("D5").Formula = '="246.16"'
$ws.Range("D6").Formula = '="1.000"'
$ws.Range("D7").Formula = '="0.4927"'
$ws.Range("D8").Formula = '="0.2957"'
$ws.Range("D9").Formula = '="0.06812"'
$ws.Range("D10").Formula = '="1.895.33"'
$ws.Range("D11").Formula = '="17.34"'
$ws.Range("D12").Formula = '="92.23"'
$ws.Range("D13").Formula = '="0.07261"'
$ws.Range("D14").Formula = '="0.6844"'
$ws.Range("D15").Formula = '="5.108"'
$ws.Range("D16").Formula = '="30.710.27"'
$ws.Range("D17").Formula = '="0.000007994"'
$ws.Range("D18").Formula = '="13.35"'
$ws.Range("D19").Formula = '="0.9999"'
$ws.Range("D20").Formula = '="2.141.01"'
$ws.Range("D21").Formula = '="1.000"'
$ws.Range("D22").Formula = '="4.864"'
$ws.Range("D23").Formula = '="190.33"'
$ws.Range("D24").Formula = '="6.085"'
$ws.Range("D25").Formula = '="9.411"'
$ws.Range("D26").Formula = '="155.63"'
$ws.Range("D27").Formula = '="19.24"'
$ws.Range("D28").Formula = '="1.933"'
$ws.Range("D29").Formula = '="1.400"'
$ws.Range("D30").Formula = '="4.396"'
$ws.Range("D31").Formula = '="0.09018"'
$ws.Range("D32").Formula = '="4.059"'
$ws.Range("D33").Formula = '="0.05202"'
$ws.Range("D34").Formula = '="0.7492"'
$ws.Range("D35").Formula = '="1.130"'
$ws.Range("D36").Formula = '="2.708"'
$ws.Range("D37").Formula = '="0.01871"'
$ws.Range("D38").Formula = '="2.675"'
$ws.Range("D39").Formula = '="2.170"'
$ws.Range("D40").Formula = '="0.9378"'
$ws.Range("D41").Formula = '="0.4452"'
$ws.Range("D42").Formula = '="106.30"'
$ws.Range("D43").Formula = '="5.826"'
$ws.Range("D44").Formula = '="1.000"'
$ws.Range("D45").Formula = '="7.731"'
$ws.Range("D46").Formula = '="0.1344"'
$ws.Range("D47").Formula = '="0.05860"'
$ws.Range("D48").Formula = '="8.801"'
$ws.Range("D49").Formula = '="0.3978"'
$ws.Range("D50").Formula = '="1.423"'
$ws.Range("D51").Formula = '="33.62"'

$dPriceRange = $ws.Range("D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51")
$dPriceRange.Copy()
$dPriceRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
